# Automatische test-sync: 2025-07-22 12:23:50
# Adds Testmail #6 (EcoPro-600 vs EcoPro-700 product-info question) to the
# "Logs" sheet, extends the conditional formatting ranges to cover the new
# row, and refreshes the category counts/order on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append the new row of mail-log data (row 7)
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A7").Value = "Waarom zit er verschil tussen de EcoPro-600 en EcoPro-700?"
$logs.Range("B7").Value = "mailmind.test@zohomail.eu"
$logs.Range("C7").Value = "Testmail #6: Waarom zit er verschil tussen de EcoPro-600 en EcoPro-700?"
$logs.Range("D7").Value = "Productinformatie"

$antwoord = "Beste klant,`r`n" +
    "Bedankt voor uw vraag over de EcoPro-600 en EcoPro-700. Het verschil tussen deze twee modellen zit voornamelijk in de capaciteit en functionaliteiten. De EcoPro-700 heeft bijvoorbeeld een grotere watertank en een extra reinigingsfunctie ten opzichte van de EcoPro-600.`r`n" +
    "Indien u meer specifieke informatie wilt over de verschillen tussen deze modellen, kunt u de productpagina’s op onze website raadplegen of contact met ons opnemen voor een gedetailleerdere uitleg.`r`n" +
    "Met vriendelijke groet,`r`n" +
    "[Naam] `r`n" +
    "E-mailassistent `r`n" +
    "--------------------------------------------------------------------------`r`n" +
    "  Dit is een testmail. Gelieve hier niet op te antwoorden."
$logs.Range("E7").Value = $antwoord

$logs.Range("F7").Value = "2025-07-22 12:23:14"
$logs.Range("G7").Value = "Ja"
$logs.Range("H7").Value = "Nee"
$logs.Range("I7").Value = "Ja"
$logs.Range("J7").Value = "Ja"

# ---------------------------------------------------------------------
# 2. Logs sheet: extend the conditional formatting sqref ranges from
#    row 6 down to row 7 (columns D, G, H, I, J)
# ---------------------------------------------------------------------
$logs.Range("D2:D6").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D7"))
$logs.Range("G2:G6").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G7"))
$logs.Range("H2:H6").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H7"))
$logs.Range("I2:I6").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I7"))
$logs.Range("J2:J6").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J7"))

# ---------------------------------------------------------------------
# 3. Dashboard sheet: refresh the category breakdown table so the new
#    "Productinformatie" mail is reflected and the rows are re-ordered
#    to match the updated counts.
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A2").Value = "Productinformatie"
$dash.Range("B2").Value = 2

$dash.Range("A3").Value = "Retour / Terugbetaling"
$dash.Range("B3").Value = 2

$dash.Range("A4").Value = "Bestelling / Levering"
$dash.Range("B4").Value = 1

$dash.Range("A5").Value = "Openingstijden / Locatie"
$dash.Range("B5").Value = 1
